$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the old row 4 ("Timepoint Type") to hold the
# new "Reference Strain" property.
$ws.Rows(4).Insert()

# Populate the new row's Property (A) and Description (C) cells.
$ws.Range("A4").Value2 = "Reference Strain"
$ws.Range("C4").Value2 = "The Reference Strain (for relative quantification data sets, leave empty for absolute)"

# Match formatting of the other rows: an empty, bold-styled B cell (like
# B2) and an italic/grey description style in C (like C3).
$ws.Range("B2").Copy()
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C4").PasteSpecial(-4122)

# Reword the "Value Unit" description (now on row 7 after the insert).
$ws.Range("C7").Value2 = "One of mM, uM, Percent, RatioT1, RatioCs, AU, Dimensionless, fmol/ug"

# Widen column A slightly to fit the new "Reference Strain" label.
$ws.Columns(1).ColumnWidth = 18.142857142857146

# Update the active cell/selection on the sheet.
$ws.Range("C18").Select() | Out-Null
